$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A : Developpeur assigné à chaque tâche
$ws.Range("A2").Value = "Anthony"
$ws.Range("A3").Value = "Anthony"
$ws.Range("A5").Value = "Anthony"
$ws.Range("A6").Value = "Olivier"
$ws.Range("A7").Value = "Anthony"
$ws.Range("A8").Value = "Anthony"
$ws.Range("A9").Value = "Olivier"
$ws.Range("A10").Value = "Olivier"
$ws.Range("A11").Value = "Olivier"
$ws.Range("A12").Value = "Olivier"

# Colonnes D (Durée estimée), E (Durée réelle), F (Difficultés)
$ws.Range("D2").Value = "30 minutes"
$ws.Range("E2").Value = "30 minutes"
$ws.Range("F2").Value = "aucune"

$ws.Range("D3").Value = "30 minutes"
$ws.Range("E3").Value = "30 minutes"
$ws.Range("F3").Value = "aucune"

$ws.Range("D5").Value = "1 heure"
$ws.Range("E5").Value = "30 minutes"
$ws.Range("F5").Value = "gerer les exceptions"

$ws.Range("D7").Value = "30 minutes"
$ws.Range("E7").Value = "1 heure"
$ws.Range("F7").Value = "recherche pour ajouter les methodes au repository"

$ws.Range("D8").Value = "2 heures"
$ws.Range("E8").Value = "2 heures"

# Mise à jour de la cellule sélectionnée
$ws.Range("F8").Select() | Out-Null
